# Update hyperparameter configuration values on the (sole) worksheet to
# reflect the "spatial clustering" run. All target cells hold their
# hyperparameter values as text (the source file stores them as inline
# strings), so each cell is forced to Text format before the new value is
# written -- this prevents the COM layer from reinterpreting a numeric-
# looking string (e.g. "200", "1.0", "1e-06") as a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "C2"  = "50"
    "C3"  = "1.0"
    "D3"  = "0.1"
    "C4"  = "200"
    "D4"  = "200"
    "F4"  = "50"
    "C5"  = "10"
    "D5"  = "10"
    "F5"  = "10"
    "D6"  = "100"
    "D8"  = "1e-06"
    "C11" = "0.03"
    "E11" = "0.1"
    "C12" = "1"
    "D12" = "0.5"
    "E12" = "0.5"
    "C13" = "250"
    "D13" = "250"
    "C14" = "4"
    "E14" = "2"
    "C16" = "2000"
    "D16" = "1000"
    "E16" = "2000"
    "F16" = "2000"
    "C17" = "5"
    "D17" = "5"
    "E17" = "20"
    "F17" = "20"
    "D18" = "4"
    "E18" = "2"
    "F18" = "2"
    "C19" = "0.01"
    "D19" = "0.01"
    "E19" = "0.1"
    "F19" = "0.1"
    "C20" = "9"
    "D20" = "9"
    "E20" = "12"
    "F20" = "4"
    "C21" = "100"
    "E21" = "100"
    "F21" = "100"
    "C22" = "10"
    "D22" = "5"
    "F22" = "5"
    "D23" = "0.1"
    "E23" = "0.05"
    "F23" = "0.05"
    "C24" = "50"
    "D24" = "10"
    "E24" = "50"
    "E25" = "0.1"
    "F25" = "0.1"
    "C26" = "[8]"
    "D26" = "[8]"
    "E26" = "[64]"
    "F26" = "[64]"
    "D27" = "400"
    "E27" = "200"
    "D28" = "14"
    "E28" = "10"
    "D29" = "4"
    "E29" = "5"
    "F29" = "4"
    "D31" = "120"
    "E31" = "150"
    "F31" = "120"
    "C33" = "100"
    "D33" = "100"
    "E33" = "1000"
    "C34" = "0.1"
    "C35" = "600"
    "D35" = "600"
    "E35" = "600"
    "F35" = "800"
    "C36" = "5"
    "E36" = "4"
    "E37" = "0.01"
}

foreach ($ref in $changes.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$ref]
}
